$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace row 2 (previously user1/password1) with the new Sauce Labs demo
# credentials. Set the password cell (B2) before the username cell (A2)
# so that the shared-string table picks up "secret_sauce" ahead of
# "standard_user", matching how the workbook was actually edited.
$ws.Range("B2").Value = "secret_sauce"
$ws.Range("A2").Value = "standard_user"

# Give the new password cell a bigger, distinct font.
$ws.Range("B2").Font.Size = 14
$ws.Range("B2").Font.Name = "Arial"

# Make row 2 a bit taller to fit the larger font.
$ws.Rows.Item(2).RowHeight = 18

# Widen both columns so the longer values are fully visible.
$ws.Columns.Item(1).ColumnWidth = 22.33
$ws.Columns.Item(2).ColumnWidth = 22.33

# Move the active selection to B11.
$null = $ws.Range("B11").Select()
